$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 261: update the Value column (G) from 1.1 to 1
$ws.Range("G261").Value = 1

# Rows 291-320: Value (G) becomes 0.5, the "Delete" flag (I) is cleared,
# and the Reference (J) switches from "Assumption for testing" to the new
# "Assumption" shared string.
for ($r = 291; $r -le 320; $r++) {
    $ws.Range("G$r").Value = 0.5
    $ws.Range("I$r").ClearContents()
    $ws.Range("J$r").Value = "Assumption"
}

# Update the current selection/scroll position to mirror the reviewer's
# view after scrolling down to the rows that were just edited.
[void]$ws.Range("A256").Select()
[void]$ws.Range("G262").Select()
